$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 22156.7
$ws.Range("J17").Value = 24396.555
$ws.Range("L17").Value = 73189.66500000001
$ws.Range("N17").Value = -73525.66500000001

$ws.Range("H80").Value = 2980.3572
$ws.Range("J80").Value = 4480
$ws.Range("L80").Value = 13440
$ws.Range("N80").Value = -15436

$ws.Range("H83").Value = 2980.3572
$ws.Range("J83").Value = 4480
$ws.Range("L83").Value = 40320
$ws.Range("N83").Value = -50304

$ws.Range("H113").Value = 5949.2607
$ws.Range("I113").Value = 4191.1665
$ws.Range("K113").Value = 4191.1665
$ws.Range("M113").Value = -937.1665000000003

$ws.Range("H137").Value = 61351.7
$ws.Range("I137").Value = 79080.52
$ws.Range("J137").Value = 3099.8572
$ws.Range("K137").Value = 237241.56
$ws.Range("L137").Value = 9299.571599999999
$ws.Range("M137").Value = -234691.56
$ws.Range("N137").Value = -14399.5716

$ws.Range("H138").Value = 2885.5469
$ws.Range("I138").Value = 1711.16
$ws.Range("J138").Value = 3638.359
$ws.Range("K138").Value = 5133.48
$ws.Range("L138").Value = 10915.077
$ws.Range("M138").Value = 6.519999999999527
$ws.Range("N138").Value = -21195.077

$ws.Range("H141").Value = 9595.25
$ws.Range("I141").Value = 9595.25
$ws.Range("K141").Value = 28785.75
$ws.Range("M141").Value = -23605.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7546.8027
$ws.Range("I32").Value = 5104.4546
$ws.Range("J32").Value = 23666.3
$ws.Range("K32").Value = 5104.4546
$ws.Range("L32").Value = 23666.3
$ws.Range("M32").Value = -4817.4546
$ws.Range("N32").Value = -24240.3

$ws.Range("H41").Value = 9262.5
$ws.Range("I41").Value = 1625
$ws.Range("J41").Value = 16900
$ws.Range("K41").Value = 1625
$ws.Range("L41").Value = 16900
$ws.Range("M41").Value = -1211
$ws.Range("N41").Value = -17728

$ws.Range("H45").Value = 8932932
$ws.Range("I45").Value = 15874969
$ws.Range("K45").Value = 15874969
$ws.Range("M45").Value = -15874592

$ws.Range("H56").Value = 30000
$ws.Range("I56").Value = 30000
$ws.Range("K56").Value = 30000
$ws.Range("M56").Value = -29258

$ws.Range("H61").Value = 4810.8335
$ws.Range("I61").Value = 4287.9473
$ws.Range("K61").Value = 4287.9473
$ws.Range("M61").Value = -4075.9473

$ws.Range("H74").Value = 21580.625
$ws.Range("I74").Value = 2360.675
$ws.Range("K74").Value = 2360.675
$ws.Range("M74").Value = -1486.675

$ws.Range("H77").Value = 21580.625
$ws.Range("I77").Value = 2360.675
$ws.Range("K77").Value = 11803.375
$ws.Range("M77").Value = -7435.375

$ws.Range("H132").Value = 25793.422
$ws.Range("I132").Value = 1395
$ws.Range("J132").Value = 63204.332
$ws.Range("K132").Value = 4185
$ws.Range("L132").Value = 189612.996
$ws.Range("M132").Value = -1655
$ws.Range("N132").Value = -194672.996

$ws.Range("H136").Value = 4810.8335
$ws.Range("I136").Value = 4287.9473
$ws.Range("K136").Value = 12863.8419
$ws.Range("M136").Value = -10313.8419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 14985.875
$ws.Range("J81").Value = 14985.875
$ws.Range("L81").Value = 14985.875
$ws.Range("N81").Value = -17107.875

$ws.Range("H84").Value = 14985.875
$ws.Range("J84").Value = 14985.875
$ws.Range("L84").Value = 44957.625
$ws.Range("N84").Value = -55565.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 17021.5
$ws.Range("J28").Value = 17021.5
$ws.Range("L28").Value = 17021.5
$ws.Range("N28").Value = -17511.5

$ws.Range("H62").Value = 5624.875
$ws.Range("I62").Value = 6366.5
$ws.Range("K62").Value = 6366.5
$ws.Range("M62").Value = -5742.5

$ws.Range("H65").Value = 5624.875
$ws.Range("I65").Value = 6366.5
$ws.Range("K65").Value = 31832.5
$ws.Range("M65").Value = -28712.5

$ws.Range("H86").Value = 7625.231
$ws.Range("I86").Value = 6192.1055
$ws.Range("J86").Value = 11515.143
$ws.Range("K86").Value = 6192.1055
$ws.Range("L86").Value = 11515.143
$ws.Range("M86").Value = -5069.1055
$ws.Range("N86").Value = -13761.143

$ws.Range("H89").Value = 7625.231
$ws.Range("I89").Value = 6192.1055
$ws.Range("J89").Value = 11515.143
$ws.Range("K89").Value = 30960.5275
$ws.Range("L89").Value = 57575.715
$ws.Range("M89").Value = -25344.5275
$ws.Range("N89").Value = -68807.715

$ws.Range("H105").Value = 1008.6923
$ws.Range("I105").Value = 716.4286
$ws.Range("K105").Value = 716.4286
$ws.Range("M105").Value = 1030.5714

$ws.Range("H116").Value = 59988
$ws.Range("J116").Value = 59988
$ws.Range("L116").Value = 59988
$ws.Range("N116").Value = -69166

$ws.Range("H130").Value = 68894.914
$ws.Range("J130").Value = 68894.914
$ws.Range("L130").Value = 68894.914
$ws.Range("N130").Value = -78934.914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 51282
$ws.Range("J5").Value = 202818.8
$ws.Range("L5").Value = 608456.3999999999
$ws.Range("N5").Value = -608680.3999999999

$ws.Range("H80").Value = 2225
$ws.Range("J80").Value = 2200
$ws.Range("L80").Value = 6600
$ws.Range("N80").Value = -8472

$ws.Range("H83").Value = 2225
$ws.Range("J83").Value = 2200
$ws.Range("L83").Value = 19800
$ws.Range("N83").Value = -29160

$ws.Range("H123").Value = 2253.375
$ws.Range("I123").Value = 2409.8
$ws.Range("J123").Value = 1992.6666
$ws.Range("K123").Value = 7229.400000000001
$ws.Range("L123").Value = 5977.9998
$ws.Range("M123").Value = -4779.400000000001
$ws.Range("N123").Value = -10877.9998

$ws.Range("H135").Value = 51282
$ws.Range("J135").Value = 202818.8
$ws.Range("L135").Value = 1825369.2
$ws.Range("N135").Value = -1830439.2

$ws.Range("H138").Value = 3231.8333
$ws.Range("J138").Value = 4579.8
$ws.Range("L138").Value = 13739.4
$ws.Range("N138").Value = -24019.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18195040
$ws.Range("I70").Value = 28574232
$ws.Range("K70").Value = 28574232
$ws.Range("M70").Value = -28573962

$ws.Range("H73").Value = 18195040
$ws.Range("I73").Value = 28574232
$ws.Range("K73").Value = 28574232
$ws.Range("M73").Value = -28573296

$ws.Range("H80").Value = 1421219.8
$ws.Range("J80").Value = 205153.81
$ws.Range("L80").Value = 205153.81
$ws.Range("N80").Value = -207149.81

$ws.Range("H83").Value = 1421219.8
$ws.Range("J83").Value = 205153.81
$ws.Range("L83").Value = 1025769.05
$ws.Range("N83").Value = -1035753.05

$ws.Range("H104").Value = 29999
$ws.Range("J104").Value = 29999
$ws.Range("L104").Value = 29999
$ws.Range("N104").Value = -36987

$ws.Range("H132").Value = 2481.0435
$ws.Range("I132").Value = 2281.75
$ws.Range("J132").Value = 2936.5715
$ws.Range("K132").Value = 6845.25
$ws.Range("L132").Value = 8809.7145
$ws.Range("M132").Value = -4315.25
$ws.Range("N132").Value = -13869.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8431.954
$ws.Range("I132").Value = 8690.619000000001
$ws.Range("K132").Value = 26071.857
$ws.Range("M132").Value = -23541.857

$ws.Range("H138").Value = 73183.86
$ws.Range("J138").Value = 73183.86
$ws.Range("L138").Value = 73183.86
$ws.Range("N138").Value = -83463.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H62").Value = 7849.2334
$ws.Range("J62").Value = 7949.5864
$ws.Range("L62").Value = 7949.5864
$ws.Range("N62").Value = -9197.5864

$ws.Range("H65").Value = 7849.2334
$ws.Range("J65").Value = 7949.5864
$ws.Range("L65").Value = 39747.932
$ws.Range("N65").Value = -45987.932

$ws.Range("H92").Value = 99999
$ws.Range("J92").Value = 99999
$ws.Range("L92").Value = 99999
$ws.Range("N92").Value = -104991

$ws.Range("H132").Value = 42127780
$ws.Range("I132").Value = 55564348
$ws.Range("J132").Value = 1818068.9
$ws.Range("K132").Value = 166693044
$ws.Range("L132").Value = 5454206.699999999
$ws.Range("M132").Value = -166690514
$ws.Range("N132").Value = -5459266.699999999
